# Add a new row for the PX_LAST / Last Price field (PR005) at the
# bottom of the FLDS_ID / FLDS_Mnemonic / FLDS_Description table on
# Sheet1 (row 45), then move the viewport/selection the way the user
# left it (scrolled down, cell I46 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A45").Value = "PR005"
$ws.Range("B45").Value = "PX_LAST"
$ws.Range("C45").Value = "Last Price"

# Scroll the window so row 22 is the first visible row, and leave the
# selection where the author left it after adding the new row.
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("I46").Select()
